# Updated cryptos list values (price + 1h volume change) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain decimal number need to be pinned to Text
# format first, otherwise Excel auto-converts the assigned string into a
# numeric value (e.g. "249.83" -> 249.8300000000001, "1.00" -> 1) and the
# trailing/leading zero formatting the site renders would be lost.
foreach ($addr in @(
    'D5',
    'D6',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D14',
    'D16',
    'D19',
    'D21',
    'D22',
    'D23',
    'D25',
    'D26',
    'D28',
    'D29',
    'D30',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D37',
    'D39',
    'D41',
    'D43',
    'D44',
    'D45',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51',
)) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '36.941.38'
$ws.Range('E2').Value = '  +4.43%  '
$ws.Range('D3').Value = '1.916.32'
$ws.Range('E3').Value = '  +1.58%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '249.83'
$ws.Range('E5').Value = '  +1.22%  '
$ws.Range('D6').Value = '0.700'
$ws.Range('E6').Value = '  +0.37%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '47.93'
$ws.Range('E8').Value = '  +10.92%  '
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  +6.26%  '
$ws.Range('D10').Value = '58.11'
$ws.Range('E10').Value = '  +6.81%  '
$ws.Range('D11').Value = '0.0761'
$ws.Range('E11').Value = '  +1.91%  '
$ws.Range('D12').Value = '0.101'
$ws.Range('E12').Value = '  +2.59%  '
$ws.Range('D13').Value = '15.14'
$ws.Range('E13').Value = '  +10.95%  '
$ws.Range('D14').Value = '0.825'
$ws.Range('E14').Value = '  +7.44%  '
$ws.Range('D15').Value = '2.194.75'
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '5.13'
$ws.Range('E16').Value = '  +3.54%  '
$ws.Range('D17').Value = '1.921.59'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('D18').Value = '36.808.75'
$ws.Range('E18').Value = '  +4.06%  '
$ws.Range('D19').Value = '74.67'
$ws.Range('E19').Value = '  +1.63%  '
$ws.Range('D20').Value = '0.0₃0858'
$ws.Range('E20').Value = '  +3.63%  '
$ws.Range('D21').Value = '13.76'
$ws.Range('E21').Value = '  +7.39%  '
$ws.Range('D22').Value = '251.93'
$ws.Range('E22').Value = '  +3.13%  '
$ws.Range('D23').Value = '5.17'
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('E24').Value = '  +0.13%  '
$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  -7.00%  '
$ws.Range('D26').Value = '2.21'
$ws.Range('E26').Value = '  +2.16%  '
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('D28').Value = '8.84'
$ws.Range('E28').Value = '  +2.34%  '
$ws.Range('D29').Value = '18.73'
$ws.Range('E29').Value = '  +2.36%  '
$ws.Range('D30').Value = '0.130'
$ws.Range('E30').Value = '  +1.12%  '
$ws.Range('E31').Value = '  +7.56%  '
$ws.Range('D32').Value = '0.0613'
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').Value = '4.33'
$ws.Range('E33').Value = '  +3.80%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.0895'
$ws.Range('E34').Value = '  +22.62%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.89'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').Value = '19.47'
$ws.Range('E37').Value = '  +60.14%  '
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('D39').Value = '0.885'
$ws.Range('E39').Value = '  +3.69%  '
$ws.Range('E40').Value = '  +1.75%  '
$ws.Range('D41').Value = '105.01'
$ws.Range('E41').Value = '  +7.89%  '
$ws.Range('E42').Value = '  +3.93%  '
$ws.Range('D43').Value = '17.76'
$ws.Range('E43').Value = '  +3.24%  '
$ws.Range('D44').Value = '2.88'
$ws.Range('E44').Value = '  +20.13%  '
$ws.Range('D45').Value = '1.10'
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('D46').Value = '1.362.03'
$ws.Range('E46').Value = '  +4.25%  '
$ws.Range('D47').Value = '2.40'
$ws.Range('E47').Value = '  -0.08%  '
$ws.Range('D48').Value = '0.0825'
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').Value = '2.81'
$ws.Range('E49').Value = '  +2.54%  '
$ws.Range('D50').Value = '6.42'
$ws.Range('E50').Value = '  +1.40%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '43.11'
$ws.Range('E51').Value = '  +2.70%  '
